$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference to its new text value. NumberFormat is forced to
# "@" (Text) before assignment so purely numeric-looking strings (for
# example "0.9975" or "1.099") are preserved as text instead of being
# coerced into floating point numbers by Excel type inference.
$updates = [ordered]@{
    "D2" = "22.140.53"
    "E2" = "  -1.00%  "
    "D3" = "1.559.96"
    "E3" = "  -0.16%  "
    "D4" = "0.9975"
    "E4" = "  -0.20%  "
    "D5" = "0.9982"
    "E5" = "  -0.15%  "
    "D6" = "292.64"
    "E6" = "  +1.65%  "
    "D7" = "0.3962"
    "E7" = "  +4.33%  "
    "D8" = "0.3258"
    "E8" = "  -0.80%  "
    "D9" = "43.80"
    "E9" = "  -1.58%  "
    "D10" = "0.07364"
    "E10" = "  -0.40%  "
    "D11" = "1.099"
    "E11" = "  -3.86%  "
    "D12" = "0.9978"
    "E12" = "  -0.18%  "
    "D13" = "19.31"
    "E13" = "  -5.71%  "
    "B14" = "Polkadot"
    "C14" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D14" = "5.685"
    "E14" = "  -3.19%  "
    "B15" = "ShibaInu"
    "C15" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "D15" = "0.00001150"
    "E15" = "  +5.72%  "
    "E16" = "  -1.52%  "
    "D17" = "1.557.73"
    "E17" = "  -0.18%  "
    "D18" = "0.06625"
    "E18" = "  -0.70%  "
    "D19" = "84.52"
    "D20" = "0.9983"
    "E20" = "  -0.08%  "
    "D21" = "6.354"
    "E21" = "  -1.51%  "
    "D22" = "15.92"
    "E22" = "  -2.02%  "
    "D23" = "11.37"
    "E23" = "  -2.87%  "
    "D24" = "22.151.21"
    "E24" = "  -0.98%  "
    "D25" = "2.345"
    "E25" = "  +2.18%  "
    "D26" = "2.482"
    "E26" = "  -4.34%  "
    "D27" = "148.77"
    "E27" = "  -0.84%  "
    "D28" = "18.75"
    "E28" = "  -3.59%  "
    "D29" = "4.879"
    "E29" = "  -1.11%  "
    "D30" = "1.732.61"
    "E30" = "  -0.43%  "
    "D31" = "119.78"
    "E31" = "  -2.36%  "
    "D32" = "1.065"
    "E32" = "  -1.56%  "
    "D33" = "5.766"
    "E33" = "  -3.45%  "
    "D34" = "0.08406"
    "E34" = "  +1.09%  "
    "D35" = "9.208"
    "E35" = "  -2.67%  "
    "D36" = "1.622"
    "E36" = "  -14.84%  "
    "D37" = "0.06227"
    "E37" = "  -1.20%  "
    "D38" = "0.02291"
    "E38" = "  -4.29%  "
    "D39" = "5.213"
    "E39" = "  -2.64%  "
    "D40" = "0.2087"
    "E40" = "  -4.02%  "
    "D41" = "1.220"
    "E41" = "  -4.63%  "
    "D42" = "10.89"
    "E42" = "  -2.23%  "
    "E43" = "  -0.20%  "
    "D44" = "0.5889"
    "E44" = "  -3.48%  "
    "D45" = "13.19"
    "E45" = "  -4.49%  "
    "D46" = "3.778"
    "E46" = "  +0.75%  "
    "D47" = "0.5646"
    "E47" = "  -5.11%  "
    "D48" = "1.921"
    "E48" = "  -4.20%  "
    "D49" = "119.14"
    "E49" = "  -4.13%  "
    "D50" = "1.149"
    "E50" = "  -2.67%  "
    "D51" = "0.06884"
    "E51" = "  -3.18%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
